$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos list (price + 1h volume change) for the daily GitHub
# Actions refresh. Column D values are leading-apostrophe-prefixed so Excel
# keeps storing them as text (matching the source inlineStr cells) instead
# of auto-coercing number-looking strings like "69.873.15" or "0.999" into
# numeric/date values.
$ws.Range("D2").Value = "'69.873.15"
$ws.Range("E2").Value = "  +2.91%  "
$ws.Range("D3").Value = "'2.578.84"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'603.52"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "'177.96"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'2.576.48"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  +17.73%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "'0.346"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "'5.05"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'0.0000185"
$ws.Range("E14").Value = "  +8.18%  "
$ws.Range("D15").Value = "'3.038.42"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "'26.36"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "'69.800.50"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "'2.563.17"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").Value = "'7.84"
$ws.Range("E19").Value = "  +3.97%  "
$ws.Range("D20").Value = "'11.24"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").Value = "'365.30"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").Value = "'4.17"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'71.02"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'4.33"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "'1.86"
$ws.Range("E26").Value = "  +6.58%  "
$ws.Range("D27").Value = "'9.36"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").Value = "'2.706.07"
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'0.0₃0929"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").Value = "'515.72"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'7.86"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'1.29"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'164.05"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'0.120"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "'19.03"
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "'18.93"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'4.95"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "'0.325"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "'2.47"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "'152.85"
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").Value = "'3.63"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "'0.526"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "'1.64"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.0₆0256"
$ws.Range("E51").Value = "  -0.68%  "
